# "Generate Report for Handback"
#
# A new handback run was generated for the 7e781ddf-d200-4b9d-b9e1-83eb1cc50d51
# doc. The latest handback file for it turned out to be stale (the version
# referenced in the handback zip isn't the tip of the branch any more), so the
# report gains an "Error Detail" entry for that row on both the zh-cn and
# de-de sheets, plus the usual Latest-Target-File / Latest-Handback-File /
# Latest-Handback-DateTime columns that a normal (non-error) row would carry.
#
# The "Error Detail" column (P) is widened to fit the long message, matching
# the existing convention used for the other long-text columns (A, G, I, J)
# in these two sheets.

$wb = $excel.ActiveWorkbook

$zhHandoffUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b360b00e9d2dc4e1c0e1aa00d517927aba79c9d0/e2e/7e781ddf-d200-4b9d-b9e1-83eb1cc50d51.md"
$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7c8be198e57325c59a15fdaded27f5d006b23f1/e2e/7e781ddf-d200-4b9d-b9e1-83eb1cc50d51.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b360b00e9d2dc4e1c0e1aa00d517927aba79c9d0/e2e/7e781ddf-d200-4b9d-b9e1-83eb1cc50d51.md."

function Update-HandbackRow {
    param(
        $ws,
        [string]$HandoffDatetime,
        [string]$TargetFile,
        [string]$HandbackDatetime
    )

    # Latest Handoff Datetime (H7) reflects the new handoff run.
    $ws.Range("H7").Value = $HandoffDatetime

    # Latest Target File (I7): same doc, displayed/linked like column A.
    $ws.Range("I7").Value = "7e781ddf-d200-4b9d-b9e1-83eb1cc50d51.md"
    $ws.Hyperlinks.Add($ws.Range("I7"), $zhHandoffUrl, "", "", "7e781ddf-d200-4b9d-b9e1-83eb1cc50d51.md") | Out-Null
    $ws.Range("I7").Style = "HyperLink"

    # Latest Handback File (J7).
    $ws.Range("J7").Value = $TargetFile

    # Latest Handback DateTime (K7).
    $ws.Range("K7").Value = $HandbackDatetime

    # Error Detail (P7): the handback version mismatch.
    $ws.Range("P7").Value = "2016-08-27 06:40:16"

    # Error Detail column needs to be wide enough for the long message.
    $ws.Columns.Item(16).ColumnWidth = 39.14
}

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow -ws $wsZh -HandoffDatetime "2016-08-27 06:40:37" -TargetFile "7e781ddf-d200-4b9d-b9e1-83eb1cc50d51.321d35a214ae4dc5b9b724be7b44a33d1cc03988.zh-cn.xlf" -HandbackDatetime $errorMessage

# de-de sheet (handoff datetime for this row is unchanged; handback datetime
# is a plain timestamp here, not an error string)
$wsDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow -ws $wsDe -HandoffDatetime "2016-08-27 06:40:20" -TargetFile "7e781ddf-d200-4b9d-b9e1-83eb1cc50d51.321d35a214ae4dc5b9b724be7b44a33d1cc03988.de-de.xlf" -HandbackDatetime "2016-08-27 06:40:43"
